$p = $ppt.ActivePresentation

function Set-ShapeText($slideIdx, $shapeIdx, $text) {
    $tr = $p.Slides.Item($slideIdx).Shapes.Item($shapeIdx).TextFrame.TextRange
    # Force a real text-range rewrite even when the final string matches the
    # existing (multi-run) concatenation, so PowerPoint collapses the runs.
    $tr.Text = "__tmp__"
    $tr.Text = $text
}

Set-ShapeText 1 1 "Slide 1 (Content)"
Set-ShapeText 2 1 "Slide 2 (Content)"
Set-ShapeText 3 1 "Slide 3 (Content)"
Set-ShapeText 4 1 "Slide 4 (Content)"
Set-ShapeText 5 1 "Slide 5 (Two Content)"

Set-ShapeText 6 1 "Slide 6 (Two Content Right)"
Set-ShapeText 6 3 "an image"

Set-ShapeText 7 1 "Slide 7 (Content with Caption)"
Set-ShapeText 7 4 "An image"

Set-ShapeText 8 1 "Slide 8 (Comparison)"
Set-ShapeText 8 4 "An image"

Set-ShapeText 9 1 "Slide 10 (Content)"
Set-ShapeText 10 1 "Slide 11 (Content)"
Set-ShapeText 11 1 "Slide 12 (Content)"
